$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.813.85"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.89%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.875.19"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.83%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.30%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.97%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.25%  "

$ws.Range("E7").Value = "  -0.67%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.12"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.98%  "

$ws.Range("E9").Value = "  +3.13%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06582"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.39%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.881.63"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.33%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.89"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.04%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07169"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.67%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6654"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.50%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "85.40"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.87%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.804"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.63%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.852.59"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.17%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007833"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.35%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.29%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.74"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.03%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.127.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.00%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9991"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.14%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.719"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.56%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.083"
$ws.Range("D24").Style = "Normal"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.554"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.07%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "147.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.26%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "134.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.15%  "

$ws.Range("E28").Value = "  +0.52%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.920"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.68%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.378"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.90%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.145"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.32%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08677"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.42%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.913"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.44%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05054"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.30%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.103"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.59%  "

$ws.Range("E36").Value = "  +3.08%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.672"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.83%  "

$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.707"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.31%  "

$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.208"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.85%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9344"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.42%  "

$ws.Range("E41").Value = "  +2.47%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.067"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.02%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9963"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.24%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "102.35"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.42%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4164"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.44%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.423"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.02%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1255"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.00%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05702"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.07%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "32.42"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.12%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.205"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.65%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "55.67"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.85%  "
